# Apply "more layout work done with added version for separate PCB box housing"
# to the sensor overview workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5 & 6: remove the stray "x" marker in column J ---------------------
$ws.Range("J5").ClearContents()
$ws.Range("J6").ClearContents()

# --- Row 20: "Relay Control" -> "Relay Control (relay box)", mark column J --
$ws.Range("A20").Value = "Relay Control (relay box)"
$ws.Range("J20").Value = "x"

# --- Row 21: H-bridge module row now needs 8 connector pins -----------------
$ws.Range("C21").Value = 8
$ws.Range("J21").Value = "x"

# --- Row 22: new "Relay Control (Expansion)" row -----------------------------
$ws.Range("A22").Value = "Relay Control (Expansion)"
$ws.Range("B22").Value = 4
$ws.Range("C22").Value = 1
$ws.Range("F22").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("H22").Value = 4
$ws.Range("J22").Value = "x"

# --- Row 23: new "Relay control power" row -----------------------------------
$ws.Range("A23").Value = "Relay control power"
$ws.Range("B23").Value = 2
$ws.Range("C23").Value = 3
$ws.Range("F23").Value = 0
$ws.Range("G23").Value = 0
$ws.Range("H23").Value = 0
$ws.Range("J23").Value = "x"

# --- Refresh the active selection to reflect where the author left off ------
$ws.Range("C21").Select()
